# Refresh the cryptos list (prices / 1h volume deltas) per the
# "Updated cryptos list ... with GitHub Actions" commit.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text still parses as a plain number need to be
# pinned to Text format first, otherwise Excel would silently convert
# them (e.g. "0.999" -> 0.999, "0.0740" -> 0.074) and we'd lose the
# exact formatting / string type that the source data uses.
$ws.Cells.Item(2, 4).Value = '43.119.63'
$ws.Cells.Item(2, 5).Value = '  +2.72%  '
$ws.Cells.Item(3, 4).Value = '2.302.01'
$ws.Cells.Item(3, 5).Value = '  +1.70%  '
$ws.Cells.Item(4, 5).Value = '  -0.12%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '310.17'
$ws.Cells.Item(5, 5).Value = '  +1.73%  '
$ws.Cells.Item(6, 5).Value = '  +5.23%  '
$ws.Cells.Item(7, 5).Value = '  +1.99%  '
$ws.Cells.Item(8, 5).Value = '  -0.03%  '
$ws.Cells.Item(9, 5).Value = '  +5.95%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '36.39'
$ws.Cells.Item(10, 5).Value = '  +3.85%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0824'
$ws.Cells.Item(11, 5).Value = '  +4.55%  '
$ws.Cells.Item(12, 5).Value = '  +0.73%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '7.14'
$ws.Cells.Item(13, 5).Value = '  +7.69%  '
$ws.Cells.Item(14, 4).Value = '2.658.18'
$ws.Cells.Item(14, 5).Value = '  +1.63%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '14.95'
$ws.Cells.Item(15, 5).Value = '  +3.86%  '
$ws.Cells.Item(16, 4).Value = '2.317.04'
$ws.Cells.Item(16, 5).Value = '  +2.14%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.808'
$ws.Cells.Item(17, 5).Value = '  +2.03%  '
$ws.Cells.Item(18, 4).Value = '43.010.54'
$ws.Cells.Item(18, 5).Value = '  +2.68%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '12.51'
$ws.Cells.Item(19, 5).Value = '  +1.19%  '
$ws.Cells.Item(20, 5).Value = '  +2.66%  '
$ws.Cells.Item(21, 5).Value = '  +1.99%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '68.35'
$ws.Cells.Item(22, 5).Value = '  +0.85%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '240.26'
$ws.Cells.Item(23, 5).Value = '  +1.17%  '
$ws.Cells.Item(24, 5).Value = '  +4.84%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.63'
$ws.Cells.Item(25, 5).Value = '  +2.31%  '
$ws.Cells.Item(26, 5).Value = '  +0.16%  '
$ws.Cells.Item(27, 5).Value = '  +3.09%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '38.56'
$ws.Cells.Item(28, 5).Value = '  +5.42%  '
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.17'
$ws.Cells.Item(29, 5).Value = '  +2.39%  '
$ws.Cells.Item(30, 2).Value = 'Cosmos'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '9.66'
$ws.Cells.Item(30, 5).Value = '  +1.84%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '167.76'
$ws.Cells.Item(31, 5).Value = '  +4.73%  '
$ws.Cells.Item(32, 5).Value = '  +2.55%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.999'
$ws.Cells.Item(33, 5).Value = '  -0.08%  '
$ws.Cells.Item(34, 5).Value = '  -0.82%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '17.80'
$ws.Cells.Item(35, 5).Value = '  +3.53%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0740'
$ws.Cells.Item(36, 5).Value = '  +0.52%  '
$ws.Cells.Item(37, 5).Value = '  +0.31%  '
$ws.Cells.Item(38, 5).Value = '  +0.79%  '
$ws.Cells.Item(39, 2).Value = 'ARBITRUM'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.84'
$ws.Cells.Item(39, 5).Value = '  +1.04%  '
$ws.Cells.Item(40, 2).Value = 'Stellar'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.116'
$ws.Cells.Item(40, 5).Value = '  +1.83%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '4.23'
$ws.Cells.Item(41, 5).Value = '  +5.73%  '
$ws.Cells.Item(42, 5).Value = '  -4.36%  '
$ws.Cells.Item(43, 2).Value = 'VeChain'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.0289'
$ws.Cells.Item(43, 5).Value = '  +2.04%  '
$ws.Cells.Item(44, 2).Value = 'Maker'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(44, 4).Value = '1.971.12'
$ws.Cells.Item(44, 5).Value = '  -0.29%  '
$ws.Cells.Item(45, 2).Value = 'EnergySwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '19.23'
$ws.Cells.Item(45, 5).Value = '  +1.59%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.04'
$ws.Cells.Item(46, 5).Value = '  +3.61%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '9.85'
$ws.Cells.Item(47, 5).Value = '  -0.23%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '55.56'
$ws.Cells.Item(48, 5).Value = '  +4.65%  '
$ws.Cells.Item(49, 5).Value = '  +16.66%  '
$ws.Cells.Item(50, 2).Value = 'Stacks'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.55'
$ws.Cells.Item(50, 5).Value = '  +2.65%  '
$ws.Cells.Item(51, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(51, 4).Value = '2.527.21'
$ws.Cells.Item(51, 5).Value = '  +1.57%  '